# Update stats for 2025-12 (row 25 in the iServ stats sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6483
$ws.Range("D25").Value = 6044921
$ws.Range("E25").Value = 932.4265000771248
$ws.Range("F25").Value = 10.04922763537599
$ws.Range("H25").Value = 26.59928718119915
